# LOQ4087.xlsx edit:
#  - Row 13 (blank label cell, holding the professor name value that had drifted
#    one row away from its "Docentes responsáveis:" label) is removed entirely,
#    shifting every subsequent row up by one. This re-aligns every label in
#    column A with its corresponding value in columns B/C.
#  - After the shift, a handful of rows' B/C value text is corrected in place.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the stray row 13 (shifts rows 14-26 up to 13-25).
$ws.Range("A13").EntireRow.Delete()

# Fix up the value cells whose text drifted out of alignment.
$ws.Range("B10:C10").Value = "8554681 - Pedro Felipe Arce Castillo"
$ws.Range("B13:C13").Value = "Semestral"
$ws.Range("B15:C15").Value = "01/01/2020"
$ws.Range("B18:C18").Value = "8554681 - Pedro Felipe Arce Castillo"
$ws.Range("B19:C19").Value = "2 provas escritas"
$ws.Range("B20:C20").Value = "Serão avaliados os conteúdos discutidos em sala e constantes da ementa do curso. A média da disciplina será a média aritmética das duas provas."
$ws.Range("B21:C21").Value = "prova escrita com conteúdo de todo o semestre"
